$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 0.6226591760299626
$wsSummary.Range("C2").Value = 0.5771495877502945
$wsSummary.Range("D2").Value = 0.9176029962546817
$wsSummary.Range("E2").Value = 0.7086044830079538
$wsSummary.Range("F2").Value = 0.8207705192629816
$wsSummary.Range("G2").Value = 0.8972462849496443
$wsSummary.Range("H2").Value = 0.7754176661195977
$wsSummary.Range("I2").Value = 490
$wsSummary.Range("J2").Value = 359
$wsSummary.Range("K2").Value = 175
$wsSummary.Range("L2").Value = 44

# ---------------------------------------------------------------------------
# Sheet: Classification Report
# ---------------------------------------------------------------------------
$wsClassification = $wb.Worksheets.Item("Classification Report")

$wsClassification.Range("B2").Value = 0.7990867579908676
$wsClassification.Range("C2").Value = 0.3277153558052435
$wsClassification.Range("D2").Value = 0.4648074369189907

$wsClassification.Range("B3").Value = 0.5771495877502945
$wsClassification.Range("C3").Value = 0.9176029962546817
$wsClassification.Range("D3").Value = 0.7086044830079538

$wsClassification.Range("B4").Value = 0.6226591760299626
$wsClassification.Range("C4").Value = 0.6226591760299626
$wsClassification.Range("D4").Value = 0.6226591760299626
$wsClassification.Range("E4").Value = 0.6226591760299626

$wsClassification.Range("B5").Value = 0.688118172870581
$wsClassification.Range("C5").Value = 0.6226591760299626
$wsClassification.Range("D5").Value = 0.5867059599634722

$wsClassification.Range("B6").Value = 0.6881181728705811
$wsClassification.Range("C6").Value = 0.6226591760299626
$wsClassification.Range("D6").Value = 0.5867059599634722

# ---------------------------------------------------------------------------
# Sheet: Confusion Matrix
# ---------------------------------------------------------------------------
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")

$wsConfusion.Range("B2").Value = 175
$wsConfusion.Range("C2").Value = 359

$wsConfusion.Range("B3").Value = 44
$wsConfusion.Range("C3").Value = 490
